$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.608.87'
$ws.Range('E2').Value = '  +2.70%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.893.17'
$ws.Range('E3').Value = '  +0.89%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9963'
$ws.Range('E4').Value = '  -0.42%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.75'
$ws.Range('E5').Value = '  +0.48%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9967'
$ws.Range('E6').Value = '  -0.37%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4928'
$ws.Range('E7').Value = '  -0.16%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2936'
$ws.Range('E8').Value = '  +1.26%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06654'
$ws.Range('E9').Value = '  +1.05%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.881.24'
$ws.Range('E10').Value = '  -0.03%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.70'
$ws.Range('E11').Value = '  -1.09%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07226'
$ws.Range('E12').Value = '  +0.73%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6751'
$ws.Range('E13').Value = '  +1.46%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.042'
$ws.Range('E14').Value = '  +5.02%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '86.85'
$ws.Range('E15').Value = '  +1.72%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.430.39'
$ws.Range('E16').Value = '  +2.09%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007895'
$ws.Range('E17').Value = '  +0.96%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9991'
$ws.Range('E18').Value = '  -0.15%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.88'
$ws.Range('E19').Value = '  +1.15%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.122.66'
$ws.Range('E20').Value = '  -0.19%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9968'
$ws.Range('E21').Value = '  -0.33%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.797'
$ws.Range('E22').Value = '  +1.60%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.945'
$ws.Range('E23').Value = '  +7.08%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.257'
$ws.Range('E24').Value = '  +1.98%  '

# Row 25
$ws.Range('B25').Value = 'BitcoinCash'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.86'
$ws.Range('E25').Value = '  +8.55%  '

# Row 26
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.13'
$ws.Range('E26').Value = '  +3.20%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.09'
$ws.Range('E27').Value = '  +2.64%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.924'
$ws.Range('E28').Value = '  +0.14%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.429'
$ws.Range('E29').Value = '  +3.64%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.256'
$ws.Range('E30').Value = '  +2.43%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08823'
$ws.Range('E31').Value = '  +1.64%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.013'
$ws.Range('E32').Value = '  +2.53%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05228'
$ws.Range('E33').Value = '  +3.36%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7317'
$ws.Range('E34').Value = '  +4.32%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.127'
$ws.Range('E35').Value = '  +2.17%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.660'
$ws.Range('E36').Value = '  -0.57%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01846'
$ws.Range('E37').Value = '  +12.86%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.706'
$ws.Range('E38').Value = '  +0.45%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.192'
$ws.Range('E39').Value = '  -0.59%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9418'
$ws.Range('E40').Value = '  +0.86%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4296'
$ws.Range('E41').Value = '  +3.32%  '

# Row 42
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.833'
$ws.Range('E42').Value = '  -3.75%  '

# Row 43
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '104.27'
$ws.Range('E43').Value = '  +1.65%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9981'
$ws.Range('E44').Value = '  +0.15%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.582'
$ws.Range('E45').Value = '  +2.04%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1306'
$ws.Range('E46').Value = '  +4.06%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05760'
$ws.Range('E47').Value = '  +1.07%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '33.02'
$ws.Range('E48').Value = '  +1.73%  '

# Row 49
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3813'
$ws.Range('E49').Value = '  +3.38%  '

# Row 50
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.341'
$ws.Range('E50').Value = '  +1.62%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.356'
$ws.Range('E51').Value = '  +1.79%  '
